$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.137.18'
$ws.Range('E2').Value = '  +0.86%  '
$ws.Range('D3').Value = '1.600.40'
$ws.Range('E3').Value = '  +0.77%  '
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').Value = '211.87'
$ws.Range('E5').Value = '  +0.54%  '
$ws.Range('E7').Value = '  +1.19%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('E9').Value = '  -0.06%  '
$ws.Range('D10').Value = '18.27'
$ws.Range('E10').Value = '  +0.02%  '
$ws.Range('D11').Value = '0.0810'
$ws.Range('E11').Value = '  +2.37%  '
$ws.Range('D12').Value = '1.821.67'
$ws.Range('E12').Value = '  +0.72%  '
$ws.Range('D13').Value = '1.592.81'
$ws.Range('E13').Value = '  +0.28%  '
$ws.Range('D14').Value = '4.03'
$ws.Range('E14').Value = '  -0.16%  '
$ws.Range('E15').Value = '  +2.50%  '
$ws.Range('D16').Value = '26.113.01'
$ws.Range('E16').Value = '  +0.79%  '
$ws.Range('D17').Value = '60.97'
$ws.Range('E17').Value = '  +1.24%  '
$ws.Range('D18').Value = '0.0₃0728'
$ws.Range('E18').Value = '  +0.30%  '
$ws.Range('E19').Value = '  -0.15%  '
$ws.Range('D20').Value = '204.44'
$ws.Range('E20').Value = '  +4.64%  '
$ws.Range('E21').Value = '  +1.76%  '
$ws.Range('D22').Value = '9.30'
$ws.Range('E22').Value = '  -1.05%  '
$ws.Range('D23').Value = '6.04'
$ws.Range('E23').Value = '  +1.36%  '
$ws.Range('E24').Value = '  +12.70%  '
$ws.Range('D25').Value = '143.29'
$ws.Range('E25').Value = '  +1.35%  '
$ws.Range('E26').Value = '  -0.06%  '
$ws.Range('E27').Value = '  -7.39%  '
$ws.Range('E28').Value = '  +0.64%  '
$ws.Range('E29').Value = '  +0.91%  '
$ws.Range('B30').Value = 'Hedera'
$ws.Range('C30').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D30').Value = '0.0479'
$ws.Range('E30').Value = '  +1.37%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').Value = '1.16'
$ws.Range('E31').Value = '  +0.24%  '
$ws.Range('E33').Value = '  -4.00%  '
$ws.Range('E34').Value = '  -0.76%  '
$ws.Range('E35').Value = '  -0.59%  '
$ws.Range('D36').Value = '1.132.12'
$ws.Range('E36').Value = '  +2.92%  '
$ws.Range('E37').Value = '  +7.60%  '
$ws.Range('B38').Value = 'PaxDollar'
$ws.Range('C38').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D38').Value = '1.00'
$ws.Range('E38').Value = '  -0.06%  '
$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D39').Value = '0.796'
$ws.Range('E39').Value = '  +1.57%  '
$ws.Range('E40').Value = '  -0.95%  '
$ws.Range('E41').Value = '  -1.93%  '
$ws.Range('D42').Value = '0.782'
$ws.Range('E42').Value = '  -2.02%  '
$ws.Range('E43').Value = '  +1.25%  '
$ws.Range('D44').Value = '1.736.33'
$ws.Range('E44').Value = '  +0.89%  '
$ws.Range('D45').Value = '92.12'
$ws.Range('E45').Value = '  -1.00%  '
$ws.Range('D46').Value = '1.50'
$ws.Range('E46').Value = '  -1.93%  '
$ws.Range('D47').Value = '54.24'
$ws.Range('E47').Value = '  +1.69%  '
$ws.Range('D48').Value = '0.0507'
$ws.Range('E48').Value = '  -0.19%  '
$ws.Range('B49').Value = 'Mantle'
$ws.Range('C49').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D49').Value = '0.407'
$ws.Range('E49').Value = '  -0.04%  '
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').Value = '0.0₇0956'
$ws.Range('E50').Value = '  -13.55%  '
$ws.Range('E51').Value = '  +0.08%  '
